$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

# --- Text content replacements (title / author / body) ---
ReplaceText "Einstein's Relativity - Redefining Space and Time" "The Mathematical Realm: Unraveling Patterns and Structures"
ReplaceText "Elijah Stone" "Clara Hudson"
ReplaceText "elijah" "clara"
ReplaceText "stone@xyz" "hudson@xyz"
ReplaceText "In the annals of scientific history, Albert Einstein's revolutionary theory of relativity stands as a towering testament to the transformative power of human intellect" "Delving into mathematics is akin to embarking on an enthralling voyage of discovery, where the interplay of numbers, patterns, and structures unfolds before our very eyes"
ReplaceText " This profound intellectual odyssey, spanning the early 20th century, unveiled a paradigm-shifting understanding of space, time, gravity, and the underlying fabric of the universe" " It is a field that captivates the mind with its precise reasoning and abstract beauty, nurturing analytical thinking and problem-solving skills that are indispensable in our ever-evolving world"
ReplaceText " Einstein's groundbreaking work challenged long-held classical notions, upending our fundamental perceptions of reality and ushering in a new era of scientific discovery" " From the ancient civilizations that deciphered numerical systems to the modern-day advancements in computer science, mathematics has left an indomitable mark on human progress and continues to shape our understanding of the universe"
ReplaceText "In his seminal 1905 paper, Einstein introduced the concept of special relativity, delving into the intricate relationship between space and time" "Mathematics is not merely a collection of abstract concepts; it finds practical applications in numerous disciplines, enriching our lives in countless ways"
ReplaceText " His groundbreaking postulates shattered the traditional view of absolute time and distance, revealing that these concepts are relative to the observer's frame of reference" " Engineers utilize mathematical principles to design and construct robust structures, while economists leverage mathematical models to analyze market trends and predict economic behavior"
ReplaceText " This radical departure from classical physics led to the iconic equation E=mc2, which elegantly encapsulates the equivalence between mass and energy" " Mathematicians themselves engage in groundbreaking research, pushing the boundaries of knowledge and expanding our comprehension of the cosmos"
ReplaceText " Einstein's special theory of relativity transformed our understanding of the universe's fundamental building blocks" " Its universality transcends cultural and linguistic barriers, connecting individuals from all corners of the globe in a shared pursuit of understanding"
ReplaceText "Einstein's intellectual journey culminated in the formulation of general relativity, a profound theory that elucidated the nature of gravity and redefined our perception of the cosmos" "The beauty of mathematics lies in its elegance and simplicity, where complex phenomena can be distilled into concise equations or formulas"
ReplaceText " General relativity conceptualizes gravity not as a force but as a curvature of spacetime caused by the presence of mass and energy" " It is a language that describes the fundamental workings of the universe, capable of expressing the laws of physics, the intricacies of biology, and the patterns of human behavior"
ReplaceText " This elegant framework revolutionized our understanding of celestial phenomena, explaining the intricate motions of planets and galaxies and opening up new avenues for exploring the vast expanse of the universe" " By delving into this realm of numbers and relationships, we uncover hidden truths and gain a deeper appreciation for the order and harmony that underpin our existence"
ReplaceText "Einstein's theory of relativity, encompassing both special and general relativity, represents a watershed moment in scientific history" "Mathematics, a captivating field of study, invites us to unravel patterns, explore structures, and uncover hidden truths"
ReplaceText " It revolutionized our understanding of space, time, gravity, and the universe's fundamental nature" " Its practical applications span a multitude of disciplines, while its aesthetic elegance captivates the mind"
ReplaceText " Einstein's groundbreaking work challenged classical notions, introducing concepts such as the relativity of simultaneity, the equivalence of mass and energy, and the curvature of spacetime" " As we delve deeper into the mathematical realm, we cultivate analytical thinking, problem-solving skills, and a profound appreciation for the order and harmony that govern our universe"
ReplaceText " His theories have had a profound impact on various scientific disciplines, astronomy, cosmology, and astrophysics, and continue to inspire and inform our quest for a deeper understanding of the cosmos" " Mathematics transcends cultural and linguistic boundaries, "

# --- Final "Summary" paragraph tail restructuring ---

$legacyRng = $d.Content
$legacyRng.Find.Execute(" Einstein's legacy as a visionary ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$legacyStart = $legacyRng.Start
$legacyEnd = $legacyRng.End

$thinkerRng = $d.Content
$thinkerRng.Find.Execute("thinker and scientific pioneer remains an enduring testament to the transformative power of intellect and the indomitable spirit of human inquiry", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$thinkerStart = $thinkerRng.Start
$thinkerEnd = $thinkerRng.End

# Replace the rightmost run first so the leftmost run's offsets stay valid.
$rThinker = $d.Range($thinkerStart, $thinkerEnd)
$rThinker.Text = "serving as a universal language that unites individuals in a shared pursuit of understanding and progress"

$rLegacy = $d.Range($legacyStart, $legacyEnd)
$rLegacy.Text = " Mathematics transcends cultural and linguistic boundaries, "

# --- Append a new empty trailing paragraph ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

# --- Fix font name across the whole document: TimesNewToman (typo) -> Times New Roman ---
$fullRng = $d.Range(0, $d.Content.End)
$fullRng.Font.Name = "Times New Roman"